$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up two pre-existing "detect_structure" (Q) values that were
#     recomputed to 0 after the stock.yaml split -----------------------
$ws.Cells.Item(61, 17).Value = 0
$ws.Cells.Item(62, 17).Value = 0

# --- Append the new weekly bars (rows 689-697) -------------------------
# Columns: A Datetime, B Open, C High, D Low, E Close, F Adj Close,
#          G Volume, H Year, I Month, J Day, K Hour, L Minute,
#          M Second, N Week, O isPivot, P two_line_structure,
#          Q detect_structure  (R "backup" is left blank for new rows,
#          same as every row was before this sync ran)
$newRows = @(
  @(689, 45474, 1808.900024414062, 1820,               1756,               1810.599975585938, 1810.599975585938, 2191358, 2024, 7, 1,  0, 0, 0, 27, 0, 0, 0),
  @(690, 45481, 1818.400024414062, 1864.699951171875,  1770,               1846.75,            1846.75,            2148274, 2024, 7, 8,  0, 0, 0, 28, 0, 0, 0),
  @(691, 45488, 1855.949951171875, 1885.949951171875,  1811.449951171875, 1841.25,            1841.25,            2552854, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 0),
  @(692, 45495, 1830.300048828125, 1838.949951171875,  1664.599975585938, 1783.199951171875, 1783.199951171875, 3990707, 2024, 7, 22, 0, 0, 0, 30, 2, 0, 0),
  @(693, 45502, 1795,               1895,               1766.550048828125, 1877.650024414062, 1877.650024414062, 3534438, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
  @(694, 45509, 1864.050048828125, 1884,               1670,               1879.199951171875, 1879.199951171875, 3058815, 2024, 8, 5,  0, 0, 0, 32, 0, 0, 0),
  @(695, 45516, 1879.199951171875, 1904.949951171875,  1785.150024414062, 1833.949951171875, 1833.949951171875, 3362034, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
  @(696, 45523, 1833,               1948.900024414062,  1832.25,            1900.800048828125, 1900.800048828125, 3099729, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
  @(697, 45530, 1909.949951171875, 2001.599975585938,  1903.75,            1965.849975585938, 1965.849975585938, 1999134, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0)
)

foreach ($row in $newRows) {
  $r = $row[0]

  $dateCell = $ws.Cells.Item($r, 1)
  $dateCell.Value = $row[1]
  $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

  $ws.Cells.Item($r, 2).Value  = $row[2]
  $ws.Cells.Item($r, 3).Value  = $row[3]
  $ws.Cells.Item($r, 4).Value  = $row[4]
  $ws.Cells.Item($r, 5).Value  = $row[5]
  $ws.Cells.Item($r, 6).Value  = $row[6]
  $ws.Cells.Item($r, 7).Value  = $row[7]
  $ws.Cells.Item($r, 8).Value  = $row[8]
  $ws.Cells.Item($r, 9).Value  = $row[9]
  $ws.Cells.Item($r, 10).Value = $row[10]
  $ws.Cells.Item($r, 11).Value = $row[11]
  $ws.Cells.Item($r, 12).Value = $row[12]
  $ws.Cells.Item($r, 13).Value = $row[13]
  $ws.Cells.Item($r, 14).Value = $row[14]
  $ws.Cells.Item($r, 15).Value = $row[15]
  $ws.Cells.Item($r, 16).Value = $row[16]
  $ws.Cells.Item($r, 17).Value = $row[17]

  # Column R ("backup") stays blank for brand-new rows. Stamp it from
  # row 688's (still blank, at this point in the script) R cell so the
  # row actually carries an empty cell node for column R, matching how
  # every other not-yet-backed-up row in this sheet is represented,
  # rather than leaving column R out of the row entirely.
  $ws.Cells.Item(688, 18).Copy($ws.Cells.Item($r, 18))
}

# --- The (until now blank) last two pre-existing rows get their
#     "backup" column populated with 0, like the rest of the sheet ----
# NOTE: done after the new rows are appended above, since they borrow
# their still-blank R688 as the blank template.
$ws.Cells.Item(687, 18).Value = 0
$ws.Cells.Item(688, 18).Value = 0
